# Fruta / hortaliza, semanal
# Insert a new daily price record at row 286 (pushing the existing
# rows 286-410 down to 287-411) and populate it with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(286).Insert()

$ws.Range("A286").Value = 4
$ws.Range("B286").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C286").Value = "Los Lagos"
$ws.Range("D286").Value = 45205
$ws.Range("E286").Value = 10
$ws.Range("F286").Value = 100112039
$ws.Range("G286").Value = "Ciboulette"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 240
$ws.Range("K286").Value = 3500
$ws.Range("L286").Value = 3500
$ws.Range("M286").Value = 3500
$ws.Range("N286").Value = "$/docena de atados"
$ws.Range("O286").Value = "Región Metropolitana"
$ws.Range("P286").Value = 1167
$ws.Range("Q286").Value = 3
$ws.Range("R286").Value = "Hortaliza"
